$wb = $excel.ActiveWorkbook

$noteText = "Using sizing capacity from pre-existing base prototype"
$seerText = "SEER Rated AC/HP"

# ---------------------------------------------------------------
# Sheet "DMo" (sheet1): add row 6 (Cap-Tons / Any / 3.5 / SEER.. / Note..)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("DMo")

# Order matters for shared-string table allocation: the note text must be
# registered first, then the "Note" header text, then the SEER text, so
# that the resulting shared-strings indices match the target workbook.
$ws1.Range("E6").Value = $noteText
$ws1.Range("E1").Value = "Note"
$ws1.Range("A6").Value = "Cap-Tons"
$ws1.Range("B6").Value = "Any"
$ws1.Range("C6").Value = 3.5
$ws1.Range("D6").Value = $seerText

$ws1.Columns.Item(4).ColumnWidth = 17.05

$ws1.Range("E6").Select()

# ---------------------------------------------------------------
# Sheet "MFm" (sheet2): add rows 38-53, CZ01..CZ16 Cap-Tons values
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("MFm")

$czLabels = @("CZ01","CZ02","CZ03","CZ04","CZ05","CZ06","CZ07","CZ08","CZ09","CZ10","CZ11","CZ12","CZ13","CZ14","CZ15","CZ16")
$vals2 = @(2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,2.5,3,3,3,3,3,3,3)

$ws2.Range("F1").Value = "Note"

for ($i = 0; $i -lt 16; $i++) {
    $r = 38 + $i
    $ws2.Cells.Item($r, 1).Value = "Cap-Tons"
    $ws2.Cells.Item($r, 2).Value = $czLabels[$i]
    $ws2.Cells.Item($r, 3).Value = $vals2[$i]
    $ws2.Cells.Item($r, 3).NumberFormat = "0.0"
    $ws2.Cells.Item($r, 4).Value = $seerText
    $ws2.Cells.Item($r, 6).Value = $noteText
}

$ws2.Range("D38:D53").Select()

# ---------------------------------------------------------------
# Sheet "SFm" (sheet3): add rows 84-99, CZ01..CZ16 Cap-Tons values
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("SFm")

$vals3 = @(4,4,4,4,4,4,4,4,4,5,5,5,5,5,5,5)

$ws3.Range("F1").Value = "Note"

for ($i = 0; $i -lt 16; $i++) {
    $r = 84 + $i
    $ws3.Cells.Item($r, 1).Value = "Cap-Tons"
    $ws3.Cells.Item($r, 2).Value = $czLabels[$i]
    $ws3.Cells.Item($r, 3).Value = $vals3[$i]
    $ws3.Cells.Item($r, 3).NumberFormat = "0.0"
    $ws3.Cells.Item($r, 4).Value = $seerText
    $ws3.Cells.Item($r, 6).Value = $noteText
}

$ws3.Range("D84:D99").Select()

# ---------------------------------------------------------------
# Make "DMo" the active sheet/tab (was "Com" before the edit)
# ---------------------------------------------------------------
$ws1.Activate()
